# Add a new worksheet named "position" after the last existing sheet
# ("unsorted"), fill it with the same (row, col, dim3, value) data that
# is already on the "3d" sheet (copied so that values / shared-string
# usage match exactly), offset so the table starts at D3 instead of A1,
# then restore the selection state of the other sheets that Excel
# touched while doing this (the previously tab-selected "missing_values"
# sheet loses tabSelected, the "3d" sheet gets a full-range selection,
# and the new "position" sheet becomes the active / tab-selected sheet
# with cell J10 selected).

$wb = $excel.ActiveWorkbook

# Source sheet holding the identical data set we need to replicate
$src = $wb.Worksheets.Item("3d")

# Insert the new sheet at the end of the workbook (after "unsorted")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "position"

# Copy the whole used range of "3d" (A1:D19) onto "position" starting at D3
$null = $src.UsedRange.Copy($newSheet.Range("D3"))

# Update the selection shown on the "3d" sheet to span its whole table
[void]$src.Activate()
$null = $src.Range("A1:D19").Select()

# Finally, activate the new sheet and select cell J10 on it, which also
# marks it as the tab-selected sheet in the saved workbook
[void]$newSheet.Activate()
$null = $newSheet.Range("J10").Select()
